$d = $word.ActiveDocument

# Locate the literal placeholder text "<<Box no>>" in the document body and
# narrow the range down to exactly that run's text.
$rng = $d.Content
$found = $rng.Find.Execute("<<Box no>>", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the '<<Box no>>' placeholder in the document."
}

# Replace that run with a MERGEFIELD construct (begin fldChar, field
# instruction, separate fldChar, cached merge result, end fldChar) — the
# same pattern already used for the other merge fields (PONumber, Date,
# Destination, ...) elsewhere in this template.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00E17C51"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:noProof/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:noProof/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD  BoxNo  \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:noProof/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:noProof/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>«BoxNo»</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:noProof/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)
